$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the numeric results (columns C, D, E for rows 2-9) ---
$ws.Range("C2").Value = -5.8419
$ws.Range("D2").Value = 0.7134
$ws.Range("E2").Value = 2.2842

$ws.Range("C3").Value = -0.6576
$ws.Range("D3").Value = 0.5346
$ws.Range("E3").Value = 1.2849

$ws.Range("C4").Value = 0.4411
$ws.Range("D4").Value = 0.4145
$ws.Range("E4").Value = 1.0097

$ws.Range("C5").Value = 0.9759
$ws.Range("D5").Value = 0.1075
$ws.Range("E5").Value = 0.2625

$ws.Range("C6").Value = 0.6758999999999999
$ws.Range("D6").Value = 0.4411
$ws.Range("E6").Value = 1.1066

$ws.Range("C7").Value = 0.416
$ws.Range("D7").Value = 0.5944
$ws.Range("E7").Value = 1.4866

$ws.Range("C8").Value = 0.0075
$ws.Range("D8").Value = 0.7879
$ws.Range("E8").Value = 1.9432

$ws.Range("C9").Value = -0.2553
$ws.Range("D9").Value = 0.9036
$ws.Range("E9").Value = 2.2414

# --- Update the heat-map style fill colors that shifted along with the
#     recomputed values. Interior.Color takes a BGR-packed integer
#     (standard VBA RGB() layout: R + G*256 + B*65536). ---
$ws.Range("D2").Interior.Color = 0xC3EACA   # fgColor -> CAEAC3
$ws.Range("D3").Interior.Color = 0x81CA81   # fgColor -> 81CA81
$ws.Range("D4").Interior.Color = 0x5FAD45   # fgColor -> 45AD5F
$ws.Range("D6").Interior.Color = 0x66B453   # fgColor -> 53B466
$ws.Range("D7").Interior.Color = 0x97D79C   # fgColor -> 9CD797
$ws.Range("D8").Interior.Color = 0xDBF3E0   # fgColor -> E0F3DB
$ws.Range("E4").Interior.Color = 0x5CA93F   # fgColor -> 3FA95C
$ws.Range("E6").Interior.Color = 0x65B352   # fgColor -> 52B365
$ws.Range("E8").Interior.Color = 0xD5F1DB   # fgColor -> DBF1D5
